$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.023758821189403534
$ws.Range("C2").Value = 0.01035654079169035
$ws.Range("D2").Value = 0.007462464272975922
$ws.Range("E2").Value = 0.005244450643658638
$ws.Range("F2").Value = 0.00004077703852090053
$ws.Range("J2").Value = 0.12743400037288666
$ws.Range("K2").Value = 1.4354861974716187
